# Update "想去人数" (interest count) values in column F across the
# "展览", "演出", "本地生活" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 2720
$ws1.Cells.Item(5, 6).Value = 936
$ws1.Cells.Item(7, 6).Value = 2331
$ws1.Cells.Item(8, 6).Value = 1844
$ws1.Cells.Item(11, 6).Value = 2489
$ws1.Cells.Item(12, 6).Value = 558
$ws1.Cells.Item(13, 6).Value = 243
$ws1.Cells.Item(14, 6).Value = 57
$ws1.Cells.Item(18, 6).Value = 9279
$ws1.Cells.Item(20, 6).Value = 7209
$ws1.Cells.Item(21, 6).Value = 11765
$ws1.Cells.Item(24, 6).Value = 235
$ws1.Cells.Item(25, 6).Value = 361
$ws1.Cells.Item(26, 6).Value = 564
$ws1.Cells.Item(27, 6).Value = 2623
$ws1.Cells.Item(30, 6).Value = 2564
$ws1.Cells.Item(31, 6).Value = 735
$ws1.Cells.Item(33, 6).Value = 4522
$ws1.Cells.Item(34, 6).Value = 939
$ws1.Cells.Item(37, 6).Value = 537

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(9, 6).Value = 1186

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(4, 6).Value = 161

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(5, 6).Value = 2720
$ws4.Cells.Item(7, 6).Value = 936
$ws4.Cells.Item(10, 6).Value = 2331
$ws4.Cells.Item(12, 6).Value = 1844
$ws4.Cells.Item(15, 6).Value = 2489
$ws4.Cells.Item(17, 6).Value = 558
$ws4.Cells.Item(18, 6).Value = 243
$ws4.Cells.Item(19, 6).Value = 57
$ws4.Cells.Item(23, 6).Value = 9279
$ws4.Cells.Item(25, 6).Value = 7209
$ws4.Cells.Item(26, 6).Value = 11765
$ws4.Cells.Item(29, 6).Value = 235
$ws4.Cells.Item(30, 6).Value = 361
$ws4.Cells.Item(32, 6).Value = 564
$ws4.Cells.Item(34, 6).Value = 2623
$ws4.Cells.Item(40, 6).Value = 4522
$ws4.Cells.Item(45, 6).Value = 537
